$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old 3x4 data grid (rows 2-13); we will repopulate as a 4x4 grid (rows 2-17)
$ws.Range("A2:T13").Clear()

# Row 2: ECs | Cxcl12 | Itgb3 | ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl12"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 198.977211
$ws.Range("H2").Value = 596.9316329999999
$ws.Range("I2").Value = 0.440791350614085
$ws.Range("J2").Value = 0.4407913506140851
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 751.8232541536289
$ws.Range("R2").Value = 6766.40928738266
$ws.Range("S2").Value = 0.1874673064093358
$ws.Range("T2").Value = 0.1874673064093358

# Row 3: ECs | Cxcl12 | Itgb3 | FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl12"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 198.977211
$ws.Range("H3").Value = 596.9316329999999
$ws.Range("I3").Value = 0.440791350614085
$ws.Range("J3").Value = 0.4407913506140851
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 862.2485094047698
$ws.Range("R3").Value = 7760.236584642928
$ws.Range("S3").Value = 0.215001869948208
$ws.Range("T3").Value = 0.2150018699482081

# Row 4: ECs | Cxcl12 | Itgb3 | M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl12"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 198.977211
$ws.Range("H4").Value = 596.9316329999999
$ws.Range("I4").Value = 0.440791350614085
$ws.Range("J4").Value = 0.4407913506140851
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 57.899450068572
$ws.Range("R4").Value = 521.095050617148
$ws.Range("S4").Value = 0.01443724158167506
$ws.Range("T4").Value = 0.01443724158167506

# Row 5: ECs | Cxcl12 | Itgb3 | sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cxcl12"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 198.977211
$ws.Range("H5").Value = 596.9316329999999
$ws.Range("I5").Value = 0.440791350614085
$ws.Range("J5").Value = 0.4407913506140851
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 95.78869058719198
$ws.Range("R5").Value = 862.0982152847279
$ws.Range("S5").Value = 0.02388493267486611
$ws.Range("T5").Value = 0.02388493267486612

# Row 6: FAPs | Cxcl12 | Itgb3 | ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl12"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 196.7746736666667
$ws.Range("H6").Value = 590.324021
$ws.Range("I6").Value = 0.4359121013721307
$ws.Range("J6").Value = 0.4359121013721308
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 743.5011011944064
$ws.Range("R6").Value = 6691.509910749657
$ws.Range("S6").Value = 0.1853921755987728
$ws.Range("T6").Value = 0.1853921755987729

# Row 7: FAPs | Cxcl12 | Itgb3 | FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl12"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 196.7746736666667
$ws.Range("H7").Value = 590.324021
$ws.Range("I7").Value = 0.4359121013721307
$ws.Range("J7").Value = 0.4359121013721308
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 852.7040267827122
$ws.Range("R7").Value = 7674.33624104441
$ws.Range("S7").Value = 0.2126219509468436
$ws.Range("T7").Value = 0.2126219509468436

# Row 8: FAPs | Cxcl12 | Itgb3 | M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Cxcl12"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 196.7746736666667
$ws.Range("H8").Value = 590.324021
$ws.Range("I8").Value = 0.4359121013721307
$ws.Range("J8").Value = 0.4359121013721308
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 57.2585440084529
$ws.Range("R8").Value = 515.3268960760761
$ws.Range("S8").Value = 0.01427743150385669
$ws.Range("T8").Value = 0.01427743150385669

# Row 9: FAPs | Cxcl12 | Itgb3 | sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Cxcl12"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 196.7746736666667
$ws.Range("H9").Value = 590.324021
$ws.Range("I9").Value = 0.4359121013721307
$ws.Range("J9").Value = 0.4359121013721308
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 94.72837736805955
$ws.Range("R9").Value = 852.555396312536
$ws.Range("S9").Value = 0.02362054332265761
$ws.Range("T9").Value = 0.02362054332265761

# Row 10: M2 | Cxcl12 | Itgb3 | ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Cxcl12"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.158424
$ws.Range("H10").Value = 0.475272
$ws.Range("I10").Value = 0.0003509544061791369
$ws.Range("J10").Value = 0.0003509544061791369
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 0.5985954201360001
$ws.Range("R10").Value = 5.387358781224
$ws.Range("S10").Value = 0.0001492599097219863
$ws.Range("T10").Value = 0.0001492599097219863

# Row 11: M2 | Cxcl12 | Itgb3 | FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Cxcl12"
$ws.Range("C11").Value = "Itgb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.158424
$ws.Range("H11").Value = 0.475272
$ws.Range("I11").Value = 0.0003509544061791369
$ws.Range("J11").Value = 0.0003509544061791369
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 0.68651508968
$ws.Range("R11").Value = 6.17863580712
$ws.Range("S11").Value = 0.0001711827001368258
$ws.Range("T11").Value = 0.0001711827001368258

# Row 12: M2 | Cxcl12 | Itgb3 | M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Cxcl12"
$ws.Range("C12").Value = "Itgb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.158424
$ws.Range("H12").Value = 0.475272
$ws.Range("I12").Value = 0.0003509544061791369
$ws.Range("J12").Value = 0.0003509544061791369
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 0.04609906044800001
$ws.Range("R12").Value = 0.4148915440320001
$ws.Range("S12").Value = 0.00001149481163616918
$ws.Range("T12").Value = 0.00001149481163616918

# Row 13: M2 | Cxcl12 | Itgb3 | sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Cxcl12"
$ws.Range("C13").Value = "Itgb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.158424
$ws.Range("H13").Value = 0.475272
$ws.Range("I13").Value = 0.0003509544061791369
$ws.Range("J13").Value = 0.0003509544061791369
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 0.07626615852800001
$ws.Range("R13").Value = 0.686395426752
$ws.Range("S13").Value = 0.00001901698468415556
$ws.Range("T13").Value = 0.00001901698468415557

# Row 14: sCs | Cxcl12 | Itgb3 | ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Cxcl12"
$ws.Range("C14").Value = "Itgb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 55.49875533333334
$ws.Range("H14").Value = 166.496266
$ws.Range("I14").Value = 0.122945593607605
$ws.Range("J14").Value = 0.122945593607605
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.778439
$ws.Range("N14").Value = 11.335317
$ws.Range("O14").Value = 0.4252971528324392
$ws.Range("P14").Value = 0.4252971528324392
$ws.Range("Q14").Value = 209.6986616029247
$ws.Range("R14").Value = 1887.287954426322
$ws.Range("S14").Value = 0.05228841091460853
$ws.Range("T14").Value = 0.05228841091460853

# Row 15: sCs | Cxcl12 | Itgb3 | FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Cxcl12"
$ws.Range("C15").Value = "Itgb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 55.49875533333334
$ws.Range("H15").Value = 166.496266
$ws.Range("I15").Value = 0.122945593607605
$ws.Range("J15").Value = 0.122945593607605
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.4877633593505858
$ws.Range("P15").Value = 0.4877633593505858
$ws.Range("Q15").Value = 240.4984913573178
$ws.Range("R15").Value = 2164.48642221586
$ws.Range("S15").Value = 0.0599683557553973
$ws.Range("T15").Value = 0.05996835575539731

# Row 16: sCs | Cxcl12 | Itgb3 | M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Cxcl12"
$ws.Range("C16").Value = "Itgb3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 55.49875533333334
$ws.Range("H16").Value = 166.496266
$ws.Range("I16").Value = 0.122945593607605
$ws.Range("J16").Value = 0.122945593607605
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2909853333333334
$ws.Range("N16").Value = 0.8729560000000001
$ws.Range("O16").Value = 0.03275300561492853
$ws.Range("P16").Value = 0.03275300561492853
$ws.Range("Q16").Value = 16.14932382025512
$ws.Range("R16").Value = 145.343914382296
$ws.Range("S16").Value = 0.004026837717760607
$ws.Range("T16").Value = 0.004026837717760607

# Row 17: sCs | Cxcl12 | Itgb3 | sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Cxcl12"
$ws.Range("C17").Value = "Itgb3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 55.49875533333334
$ws.Range("H17").Value = 166.496266
$ws.Range("I17").Value = 0.122945593607605
$ws.Range("J17").Value = 0.122945593607605
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4814053333333333
$ws.Range("N17").Value = 1.444216
$ws.Range("O17").Value = 0.0541864822020464
$ws.Range("P17").Value = 0.05418648220204641
$ws.Range("Q17").Value = 26.71739681082845
$ws.Range("R17").Value = 240.456571297456
$ws.Range("S17").Value = 0.006661989219838516
$ws.Range("T17").Value = 0.006661989219838516
